# Apply edits described by the diff:
# - Rename Sheet1 (currently "ID_0e9fca7") -> "ID_9c650ad"
# - Rename Sheet2 (currently "ID_9c650ad") -> "ID_ae0513b"
# - Update Sheet1 data row 2: price 22.3 -> 55.99, product -> "Mens Cotton Jacket"
# - Update Sheet2 data row 2: price 55.99 -> 15.99, date 04/02/2025 -> 12/02/2025,
#   product "Mens Cotton Jacket" -> "Mens Casual Slim Fit"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update data first (while names still make sheets easy to identify by position)
$ws1.Range("A2").Value = 55.99
$ws1.Range("C2").Value = "Mens Cotton Jacket"

$ws2.Range("A2").Value = 15.99
# Force the date-like text to stay plain text (matches the source inlineStr
# "12/02/2025") instead of letting Excel auto-convert it into a date serial,
# the same way typing a leading apostrophe in the Excel UI would.
$ws2.Range("B2").Value = "'12/02/2025"
$ws2.Range("C2").Value = "Mens Casual Slim Fit"

# Rename sheets. Sheet2 currently holds the name "ID_9c650ad" that Sheet1
# needs to take, so rename Sheet2 out of the way first to avoid a collision.
$ws2.Name = "ID_ae0513b"
$ws1.Name = "ID_9c650ad"
